# The deck's applied design ("Integral") is switched over to the built-in
# "Office Theme" colour palette. In the underlying package this shows up as
# the theme part that backs the slide master (ppt/theme/theme2.xml, the
# presentation's active theme) taking on the Office Theme's 12 colours
# (and, symmetrically, the palette that used to live there ends up parked in
# the other theme part used by the notes master). We drive this the same way
# a user would from the Design tab: by setting each of the twelve theme
# colour slots on the active presentation theme.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$themeColors = $master.Theme.ThemeColorScheme

function Set-ThemeRGB {
    param($Scheme, [int]$Index, [string]$Hex)

    $r = [Convert]::ToInt32($Hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($Hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($Hex.Substring(4, 2), 16)

    $Scheme.Colors($Index).RGB = $r + ($g * 256) + ($b * 65536)
}

# Office Theme colour scheme, in ThemeColorScheme.Colors() slot order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
# 9 accent5, 10 accent6, 11 hlink, 12 folHlink.
Set-ThemeRGB $themeColors 1  "000000"
Set-ThemeRGB $themeColors 2  "FFFFFF"
Set-ThemeRGB $themeColors 3  "44546A"
Set-ThemeRGB $themeColors 4  "E7E6E6"
Set-ThemeRGB $themeColors 5  "5B9BD5"
Set-ThemeRGB $themeColors 6  "ED7D31"
Set-ThemeRGB $themeColors 7  "A5A5A5"
Set-ThemeRGB $themeColors 8  "FFC000"
Set-ThemeRGB $themeColors 9  "4472C4"
Set-ThemeRGB $themeColors 10 "70AD47"
Set-ThemeRGB $themeColors 11 "0563C1"
Set-ThemeRGB $themeColors 12 "954F72"
